# Commit: "update script docstrings, change eg_dict to p_dict"
#
# Rename the "eg_dictionary" settings sheet to "proposal_dictionary" and
# update its header cell (A1) from "eg" to "proposal" to match. Also make
# this sheet the active/selected sheet (moving the selection highlight off
# of "ret_incr", which was previously the active sheet), with the cursor
# left on E12.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("eg_dictionary")
$ws.Name = "proposal_dictionary"
$ws.Range("A1").Value = "proposal"

$ws.Activate()
$ws.Range("E12").Select()
